# Update simulation results for the 380 kV case (pl_mw.xlsx, "res_line" sheet data)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("B2").Value = 0.5622756085975311
$ws.Range("D2").Value = 0.02960867911500742
$ws.Range("E2").Value = 0.3621951315892704
$ws.Range("F2").Value = 0.5475057961777239
$ws.Range("G2").Value = 0.3883662425929231
$ws.Range("H2").Value = 0.5568190128946853
$ws.Range("I2").Value = 0.9995892122024497
$ws.Range("K2").Value = 0.5202082050356296
$ws.Range("L2").Value = 0.1282284262681941
$ws.Range("M2").Value = 0.1427426126049625
$ws.Range("O2").Value = 1.834042247836109

# Row 3
$ws.Range("B3").Value = 0.5388885829726746
$ws.Range("D3").Value = 0.02667756787597142
$ws.Range("E3").Value = 0.3654509267136796
$ws.Range("F3").Value = 0.5477684406437007
$ws.Range("G3").Value = 0.3900127814993795
$ws.Range("H3").Value = 0.5609305696444551
$ws.Range("I3").Value = 1.012007347558423
$ws.Range("K3").Value = 0.4690243635214983
$ws.Range("L3").Value = 0.121790384033801
$ws.Range("M3").Value = 0.1363775074497156
$ws.Range("O3").Value = 1.845911253040782

# Row 4
$ws.Range("B4").Value = 0.5246913574374901
$ws.Range("D4").Value = 0.02486550517321717
$ws.Range("E4").Value = 0.3675681809030094
$ws.Range("F4").Value = 0.5482881788817622
$ws.Range("G4").Value = 0.3913151934372294
$ws.Range("H4").Value = 0.5637031371758567
$ws.Range("I4").Value = 1.020093972564506
$ws.Range("K4").Value = 0.4374107587989897
$ws.Range("L4").Value = 0.1178859597127655
$ws.Range("M4").Value = 0.1325153555477598
$ws.Range("O4").Value = 1.854327943518385

# Row 5
$ws.Range("B5").Value = 0.51894728618268
$ws.Range("D5").Value = 0.0241240093488102
$ws.Range("E5").Value = 0.3684607490171328
$ws.Range("F5").Value = 0.5485901548464298
$ws.Range("G5").Value = 0.3919191766058958
$ws.Range("H5").Value = 0.5648954088991687
$ws.Range("I5").Value = 1.023505515327452
$ws.Range("K5").Value = 0.4244819447153247
$ws.Range("L5").Value = 0.1163071962756277
$ws.Range("M5").Value = 0.1309532130222983
$ws.Range("O5").Value = 1.85804174151265

# Row 6
$ws.Range("B6").Value = 0.5179960050187447
$ws.Range("D6").Value = 0.02400070054394376
$ws.Range("E6").Value = 0.3686107588211121
$ws.Range("F6").Value = 0.5486457453676223
$ws.Range("G6").Value = 0.3920238897858823
$ws.Range("H6").Value = 0.5650971567471927
$ws.Range("I6").Value = 1.024079017529452
$ws.Range("K6").Value = 0.4223323720775625
$ws.Range("L6").Value = 0.1160457912293609
$ws.Range("M6").Value = 0.1306945320661264
$ws.Range("O6").Value = 1.858675563797775

# Row 7
$ws.Range("B7").Value = 0.5246137224707184
$ws.Range("D7").Value = 0.02485551746110559
$ws.Range("E7").Value = 0.3675800977823656
$ws.Range("F7").Value = 0.5482918862569903
$ws.Range("G7").Value = 0.391323042478021
$ws.Range("H7").Value = 0.5637189637313398
$ws.Range("I7").Value = 1.020139511365988
$ws.Range("K7").Value = 0.4372365814184036
$ws.Range("L7").Value = 0.1178646179225993
$ws.Range("M7").Value = 0.1324942403295672
$ws.Range("O7").Value = 1.854376879464354

# Row 8
$ws.Range("B8").Value = 0.5541783567306027
$ws.Range("D8").Value = 0.02860061783994894
$ws.Range("E8").Value = 0.3632932322709403
$ws.Range("F8").Value = 0.5475219757770304
$ws.Range("G8").Value = 0.3888734564273122
$ws.Range("H8").Value = 0.5581852252488346
$ws.Range("I8").Value = 1.003775158141909
$ws.Range("K8").Value = 0.5025992456992299
$ws.Range("L8").Value = 0.1259985748210397
$ws.Range("M8").Value = 0.1405384466419086
$ws.Range("O8").Value = 1.837900348161881

# Row 9
$ws.Range("B9").Value = 0.6134229176605572
$ws.Range("D9").Value = 0.03584534859020039
$ws.Range("E9").Value = 0.3558220341365015
$ws.Range("F9").Value = 0.5488556494041958
$ws.Range("G9").Value = 0.3863845531193775
$ws.Range("H9").Value = 0.5492995640650804
$ws.Range("I9").Value = 0.9753468533030283
$ws.Range("K9").Value = 0.6292620930045985
$ws.Range("L9").Value = 0.1423305709971743
$ws.Range("M9").Value = 0.1566733622428274
$ws.Range("O9").Value = 1.814549206055261

# Row 10
$ws.Range("B10").Value = 0.6576996767951755
$ws.Range("D10").Value = 0.04110601843690631
$ws.Range("E10").Value = 0.3508997509826581
$ws.Range("F10").Value = 0.5515689512053328
$ws.Range("G10").Value = 0.3859709550311194
$ws.Range("H10").Value = 0.5439669205611821
$ws.Range("I10").Value = 0.9566893522558235
$ws.Range("K10").Value = 0.7213628056065602
$ws.Range("L10").Value = 0.1545581308636059
$ws.Range("M10").Value = 0.1687417073076034
$ws.Range("O10").Value = 1.802857518382964

# Row 11
$ws.Range("B11").Value = 0.6780004482763786
$ws.Range("D11").Value = 0.043485484317074
$ws.Range("E11").Value = 0.348782828657942
$ws.Range("F11").Value = 0.5531797663869753
$ws.Range("G11").Value = 0.3860909127201637
$ws.Range("H11").Value = 0.541799994688354
$ws.Range("I11").Value = 0.9486847119265249
$ws.Range("K11").Value = 0.7630465840919669
$ws.Range("L11").Value = 0.1601695907241378
$ws.Range("M11").Value = 0.1742772746138712
$ws.Range("O11").Value = 1.798725866153916

# Row 12
$ws.Range("B12").Value = 0.6857101982737106
$ws.Range("D12").Value = 0.04438452996946296
$ws.Range("E12").Value = 0.3479987311866866
$ws.Range("F12").Value = 0.5538438635802621
$ws.Range("G12").Value = 0.3861806995337957
$ws.Range("H12").Value = 0.5410166208038589
$ws.Range("I12").Value = 0.9457229370459768
$ws.Range("K12").Value = 0.7787996871902862
$ws.Range("L12").Value = 0.1623014672509697
$ws.Range("M12").Value = 0.1763798879669451
$ws.Range("O12").Value = 1.797332025348254

# Row 13
$ws.Range("B13").Value = 0.684048783654589
$ws.Range("D13").Value = 0.04419099435845908
$ws.Range("E13").Value = 0.3481668214301754
$ws.Range("F13").Value = 0.5536984319664313
$ws.Range("G13").Value = 0.3861593884683003
$ws.Range("H13").Value = 0.5411836807911925
$ws.Range("I13").Value = 0.9463577214512267
$ws.Range("K13").Value = 0.7754083971051386
$ws.Range("L13").Value = 0.1618420228101343
$ws.Range("M13").Value = 0.1759267696047218
$ws.Range("O13").Value = 1.797624620148639

# Row 14
$ws.Range("B14").Value = 0.6786342906908942
$ws.Range("D14").Value = 0.04355948981508106
$ws.Range("E14").Value = 0.3487179693102345
$ws.Range("F14").Value = 0.5532333176471838
$ws.Range("G14").Value = 0.3860974101484373
$ws.Range("H14").Value = 0.5417348008472445
$ws.Range("I14").Value = 0.9484396539291247
$ws.Range("K14").Value = 0.764343240989632
$ws.Range("L14").Value = 0.1603448431166896
$ws.Range("M14").Value = 0.1744501301635282
$ws.Range("O14").Value = 1.798607772185335

# Row 15
$ws.Range("B15").Value = 0.6753206427365512
$ws.Range("D15").Value = 0.04317241240981673
$ws.Range("E15").Value = 0.3490578455781339
$ws.Range("F15").Value = 0.5529554684354494
$ws.Range("G15").Value = 0.386065225453784
$ws.Range("H15").Value = 0.5420772203959388
$ws.Range("I15").Value = 0.9497239364455847
$ws.Range("K15").Value = 0.7575613643656425
$ws.Range("L15").Value = 0.1594286773702152
$ws.Range("M15").Value = 0.1735464768652974
$ws.Range("O15").Value = 1.799232216032351

# Row 16
$ws.Range("B16").Value = 0.6563761243001238
$ws.Range("D16").Value = 0.04095023643372997
$ws.Range("E16").Value = 0.3510405530135565
$ws.Range("F16").Value = 0.551471254948531
$ws.Range("G16").Value = 0.3859693202221308
$ws.Range("H16").Value = 0.5441137411113885
$ws.Range("I16").Value = 0.9572221876209106
$ws.Range("K16").Value = 0.7186342978032201
$ws.Range("L16").Value = 0.1541923862080807
$ws.Range("M16").Value = 0.168380850451598
$ws.Range("O16").Value = 1.803151419931098

# Row 17
$ws.Range("B17").Value = 0.6447945976240987
$ws.Range("D17").Value = 0.03958347803580864
$ws.Range("E17").Value = 0.352288159170649
$ws.Range("F17").Value = 0.5506571554356654
$ws.Range("G17").Value = 0.3859894391155407
$ws.Range("H17").Value = 0.5454293689355936
$ws.Range("I17").Value = 0.9619457584831679
$ws.Range("K17").Value = 0.694698483453692
$ws.Range("L17").Value = 0.1509925772407996
$ws.Range("M17").Value = 0.1652234859358259
$ws.Range("O17").Value = 1.805859764143449

# Row 18
$ws.Range("B18").Value = 0.6381482152781643
$ws.Range("D18").Value = 0.03879607367017002
$ws.Range("E18").Value = 0.3530172585108995
$ws.Range("F18").Value = 0.5502243485217306
$ws.Range("G18").Value = 0.3860300080422903
$ws.Range("H18").Value = 0.546210456381985
$ws.Range("I18").Value = 0.9647080708717795
$ws.Range("K18").Value = 0.6809112171449101
$ws.Range("L18").Value = 0.1491567586915323
$ws.Range("M18").Value = 0.1634117557149359
$ws.Range("O18").Value = 1.807529251515675

# Row 19
$ws.Range("B19").Value = 0.6359004594040414
$ws.Range("D19").Value = 0.03852925342761182
$ws.Range("E19").Value = 0.3532660970998167
$ws.Range("F19").Value = 0.5500838954563108
$ws.Range("G19").Value = 0.3860487225239453
$ws.Range("H19").Value = 0.54647910645118
$ws.Range("I19").Value = 0.9656511469530038
$ws.Range("K19").Value = 0.6762396769173336
$ws.Range("L19").Value = 0.1485359803905908
$ws.Range("M19").Value = 0.1627990783568194
$ws.Range("O19").Value = 1.808113697799797

# Row 20
$ws.Range("B20").Value = 0.6460259221263129
$ws.Range("D20").Value = 0.03972910471406976
$ws.Range("E20").Value = 0.3521541585198433
$ws.Range("F20").Value = 0.5507401498360309
$ws.Range("G20").Value = 0.3859842959255104
$ws.Range("H20").Value = 0.5452867958604273
$ws.Range("I20").Value = 0.961438223572598
$ws.Range("K20").Value = 0.6972485706578482
$ws.Range("L20").Value = 0.1513327243618932
$ws.Range("M20").Value = 0.1655591485567811
$ws.Range("O20").Value = 1.805559893414028

# Row 21
$ws.Range("B21").Value = 0.6802240582640593
$ws.Range("D21").Value = 0.04374503282160447
$ws.Range("E21").Value = 0.3485556082718497
$ws.Range("F21").Value = 0.5533684645132979
$ws.Range("G21").Value = 0.3861144102848897
$ws.Range("H21").Value = 0.5415719144019278
$ws.Range("I21").Value = 0.9478262562425375
$ws.Range("K21").Value = 0.7675942119340675
$ws.Range("L21").Value = 0.1607844137297292
$ws.Range("M21").Value = 0.1748836820422568
$ws.Range("O21").Value = 1.79831436280702

# Row 22
$ws.Range("B22").Value = 0.7027041613337417
$ws.Range("D22").Value = 0.04635795793657849
$ws.Range("E22").Value = 0.3463059351647759
$ws.Range("F22").Value = 0.5554016386578482
$ws.Range("G22").Value = 0.3864580451349013
$ws.Range("H22").Value = 0.5393608013135349
$ws.Range("I22").Value = 0.9393346110319403
$ws.Range("K22").Value = 0.8133842912290277
$ws.Range("L22").Value = 0.1670020262435656
$ws.Range("M22").Value = 0.1810151228608348
$ws.Range("O22").Value = 1.794574130936979

# Row 23
$ws.Range("B23").Value = 0.6906944371536952
$ws.Range("D23").Value = 0.04496447758519651
$ws.Range("E23").Value = 0.3474972924009745
$ws.Range("F23").Value = 0.5542876437618602
$ws.Range("G23").Value = 0.3862509608073879
$ws.Range("H23").Value = 0.5405210916394196
$ws.Range("I23").Value = 0.9438297477537656
$ws.Range("K23").Value = 0.7889624979437144
$ws.Range("L23").Value = 0.1636799135026337
$ws.Range("M23").Value = 0.1777392908321929
$ws.Range("O23").Value = 1.796479293958555

# Row 24
$ws.Range("B24").Value = 0.6454692029231239
$ws.Range("D24").Value = 0.03966327198799036
$ws.Range("E24").Value = 0.3522147033497576
$ws.Range("F24").Value = 0.5507025183285847
$ws.Range("G24").Value = 0.3859865308264858
$ws.Range("H24").Value = 0.5453511761577658
$ws.Range("I24").Value = 0.9616675347368915
$ws.Range("K24").Value = 0.6960957579191813
$ws.Range("L24").Value = 0.1511789320211818
$ws.Range("M24").Value = 0.1654073846331414
$ws.Range("O24").Value = 1.805695114782537

# Row 25
$ws.Range("B25").Value = 0.5972623341455403
$ws.Range("D25").Value = 0.03389624663306989
$ws.Range("E25").Value = 0.3577434097383219
$ws.Range("F25").Value = 0.548190437985582
$ws.Range("G25").Value = 0.3868096546580873
$ws.Range("H25").Value = 0.5514931950317745
$ws.Range("I25").Value = 0.982645795836655
$ws.Range("K25").Value = 0.5951621663460855
$ws.Range("L25").Value = 0.1378719250622424
$ws.Range("M25").Value = 0.1522704517802573
$ws.Range("O25").Value = 1.819906880378241

